$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 365
$ws.Range("I28").Value = 365
$ws.Range("K28").Value = 365
$ws.Range("M28").Value = 120
$ws.Range("H80").Value = 417.25
$ws.Range("I80").Value = 498.1111
$ws.Range("J80").Value = 351.0909
$ws.Range("K80").Value = 1494.3333
$ws.Range("L80").Value = 1053.2727
$ws.Range("M80").Value = -496.3333
$ws.Range("N80").Value = -3049.2727
$ws.Range("H83").Value = 417.25
$ws.Range("I83").Value = 498.1111
$ws.Range("J83").Value = 351.0909
$ws.Range("K83").Value = 4482.9999
$ws.Range("L83").Value = 3159.8181
$ws.Range("M83").Value = 509.0001000000002
$ws.Range("N83").Value = -13143.8181
$ws.Range("H100").Value = 4500.1377
$ws.Range("I100").Value = 1192.3334
$ws.Range("K100").Value = 1192.3334
$ws.Range("M100").Value = -651.3334
$ws.Range("H107").Value = 430.89474
$ws.Range("I107").Value = 414.1111
$ws.Range("J107").Value = 446
$ws.Range("K107").Value = 414.1111
$ws.Range("L107").Value = 446
$ws.Range("M107").Value = 1505.8889
$ws.Range("N107").Value = -4286
$ws.Range("H129").Value = 3427
$ws.Range("I129").Value = 3230.375
$ws.Range("K129").Value = 9691.125
$ws.Range("M129").Value = -4691.125
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("H138").Value = 4405.2915
$ws.Range("I138").Value = 3611.125
$ws.Range("J138").Value = 4564.125
$ws.Range("K138").Value = 10833.375
$ws.Range("L138").Value = 13692.375
$ws.Range("M138").Value = -5693.375
$ws.Range("N138").Value = -23972.375

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2062.081
$ws.Range("I2").Value = 1707.1364
$ws.Range("J2").Value = 2582.6667
$ws.Range("K2").Value = 1707.1364
$ws.Range("L2").Value = 2582.6667
$ws.Range("M2").Value = -1594.1364
$ws.Range("N2").Value = -2808.6667
$ws.Range("H32").Value = 2390517.2
$ws.Range("I32").Value = 2737336
$ws.Range("K32").Value = 2737336
$ws.Range("M32").Value = -2737049
$ws.Range("H35").Value = 3330.3333
$ws.Range("I35").Value = 3330.3333
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3330.3333
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2924.3333
$ws.Range("N35").ClearContents()
$ws.Range("H36").Value = 50000
$ws.Range("J36").Value = 50000
$ws.Range("L36").Value = 50000
$ws.Range("N36").Value = -50692
$ws.Range("H53").Value = 9234.4
$ws.Range("J53").Value = 9234.4
$ws.Range("L53").Value = 9234.4
$ws.Range("N53").Value = -10598.4
$ws.Range("H74").Value = 3679379.8
$ws.Range("I74").Value = 4631058.5
$ws.Range("K74").Value = 4631058.5
$ws.Range("M74").Value = -4630184.5
$ws.Range("H77").Value = 3679379.8
$ws.Range("I77").Value = 4631058.5
$ws.Range("K77").Value = 23155292.5
$ws.Range("M77").Value = -23150924.5
$ws.Range("H102").Value = 2356.1155
$ws.Range("I102").Value = 2431.1667
$ws.Range("K102").Value = 2431.1667
$ws.Range("M102").Value = -809.1667000000002
$ws.Range("H116").Value = 2062.081
$ws.Range("I116").Value = 1707.1364
$ws.Range("J116").Value = 2582.6667
$ws.Range("K116").Value = 1707.1364
$ws.Range("L116").Value = 2582.6667
$ws.Range("M116").Value = 586.8635999999999
$ws.Range("N116").Value = -7170.6667

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2062.081
$ws.Range("I3").Value = 1707.1364
$ws.Range("J3").Value = 2582.6667
$ws.Range("K3").Value = 1707.1364
$ws.Range("L3").Value = 2582.6667
$ws.Range("M3").Value = -1593.1364
$ws.Range("N3").Value = -2810.6667
$ws.Range("H20").Value = 1740.7142
$ws.Range("I20").Value = 2039.7333
$ws.Range("J20").Value = 1395.6923
$ws.Range("K20").Value = 2039.7333
$ws.Range("L20").Value = 1395.6923
$ws.Range("M20").Value = -1792.7333
$ws.Range("N20").Value = -1889.6923
$ws.Range("H86").Value = 1859.625
$ws.Range("I86").Value = 1979.8334
$ws.Range("K86").Value = 1979.8334
$ws.Range("M86").Value = -856.8334
$ws.Range("H89").Value = 1859.625
$ws.Range("I89").Value = 1979.8334
$ws.Range("K89").Value = 9899.166999999999
$ws.Range("M89").Value = -4283.166999999999
$ws.Range("H94").Value = 1195.4615
$ws.Range("I94").Value = 1260.3334
$ws.Range("K94").Value = 1260.3334
$ws.Range("M94").Value = -809.3334
$ws.Range("H102").Value = 95000
$ws.Range("I102").Value = 130000
$ws.Range("K102").Value = 130000
$ws.Range("M102").Value = -126755
$ws.Range("H105").Value = 2785.318
$ws.Range("I105").Value = 2383.158
$ws.Range("J105").Value = 5332.3335
$ws.Range("K105").Value = 2383.158
$ws.Range("L105").Value = 5332.3335
$ws.Range("M105").Value = -636.1579999999999
$ws.Range("N105").Value = -8826.333500000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6691.6665
$ws.Range("I31").Value = 1924.3334
$ws.Range("K31").Value = 1924.3334
$ws.Range("M31").Value = -1629.3334
$ws.Range("H34").Value = 6691.6665
$ws.Range("I34").Value = 1924.3334
$ws.Range("K34").Value = 1924.3334
$ws.Range("M34").Value = -1722.3334

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3022.5
$ws.Range("J5").Value = 9279.6
$ws.Range("L5").Value = 27838.8
$ws.Range("N5").Value = -28062.8
$ws.Range("H8").Value = 471.25
$ws.Range("I8").Value = 471.25
$ws.Range("K8").Value = 1413.75
$ws.Range("M8").Value = -1274.75
$ws.Range("H11").Value = 4546261.5
$ws.Range("I11").Value = 314
$ws.Range("J11").Value = 6667703.5
$ws.Range("K11").Value = 942
$ws.Range("L11").Value = 20003110.5
$ws.Range("M11").Value = -802
$ws.Range("N11").Value = -20003390.5
$ws.Range("H14").Value = 133781.53
$ws.Range("I14").Value = 133781.53
$ws.Range("K14").Value = 401344.59
$ws.Range("M14").Value = -401171.59
$ws.Range("H37").Value = 99998
$ws.Range("J37").Value = 99998
$ws.Range("L37").Value = 299994
$ws.Range("N37").Value = -300218
$ws.Range("H68").Value = 11101.5
$ws.Range("J68").Value = 13126.875
$ws.Range("L68").Value = 39380.625
$ws.Range("N68").Value = -41002.625
$ws.Range("H71").Value = 11101.5
$ws.Range("J71").Value = 13126.875
$ws.Range("L71").Value = 118141.875
$ws.Range("N71").Value = -126253.875
$ws.Range("H75").Value = 8102.4546
$ws.Range("I75").Value = 4950
$ws.Range("J75").Value = 8803
$ws.Range("K75").Value = 14850
$ws.Range("L75").Value = 26409
$ws.Range("M75").Value = -13852
$ws.Range("N75").Value = -28405
$ws.Range("H78").Value = 8102.4546
$ws.Range("I78").Value = 4950
$ws.Range("J78").Value = 8803
$ws.Range("K78").Value = 44550
$ws.Range("L78").Value = 79227
$ws.Range("M78").Value = -39558
$ws.Range("N78").Value = -89211
$ws.Range("H131").Value = 8338.471
$ws.Range("J131").Value = 9870.893
$ws.Range("L131").Value = 29612.679
$ws.Range("N131").Value = -39692.679
$ws.Range("H135").Value = 3022.5
$ws.Range("J135").Value = 9279.6
$ws.Range("L135").Value = 83516.40000000001
$ws.Range("N135").Value = -88586.40000000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 390308.16
$ws.Range("J80").Value = 6039
$ws.Range("L80").Value = 6039
$ws.Range("N80").Value = -8035
$ws.Range("H83").Value = 390308.16
$ws.Range("J83").Value = 6039
$ws.Range("L83").Value = 30195
$ws.Range("N83").Value = -40179
$ws.Range("H102").Value = 3087.4211
$ws.Range("I102").Value = 2990.7334
$ws.Range("K102").Value = 2990.7334
$ws.Range("M102").Value = -1368.7334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1508157.8
$ws.Range("I132").Value = 1823559.4
$ws.Range("K132").Value = 5470678.199999999
$ws.Range("M132").Value = -5468148.199999999
$ws.Range("H136").Value = 4715.2
$ws.Range("I136").Value = 3069.25
$ws.Range("K136").Value = 9207.75
$ws.Range("M136").Value = -6657.75
$ws.Range("H137").Value = 115286
$ws.Range("J137").Value = 115286
$ws.Range("L137").Value = 115286
$ws.Range("N137").Value = -125486

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2344.4443
$ws.Range("I81").Value = 2344.4443
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4688.8886
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3627.8886
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2344.4443
$ws.Range("I84").Value = 2344.4443
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 23444.443
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -18140.443
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 9035838
$ws.Range("I132").Value = 641551.75
$ws.Range("J132").Value = 59401556
$ws.Range("K132").Value = 1924655.25
$ws.Range("L132").Value = 178204668
$ws.Range("M132").Value = -1922125.25
$ws.Range("N132").Value = -178209728
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
